# presentation: fix wording on diverted replica threshold slide.
#
# Slide 34 ("Replica Diversion") content placeholder:
#   - turn on "shrink text on overflow" autofit (<a:normAutofit/>) so the
#     extra bullet below still fits the placeholder
#   - append a new bullet asking what happens when the node picked for a
#     diverted replica can't actually store the file

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(34)
$shape = $s.Shapes.Item(3)

# msoAutoSizeTextToFitShape -> <a:bodyPr><a:normAutofit/></a:bodyPr>
$shape.TextFrame.AutoSize = 2

$tr = $shape.TextFrame.TextRange
$tr.InsertAfter("`rWhat happens when node picked for diverted replica can’t store the file?") | Out-Null
